$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "2h 28m"
$ws.Range("C2").Select() | Out-Null
